$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1389.0471
$ws.Range("I15").Value = 1389.0471
$ws.Range("K15").Value = 4167.1413
$ws.Range("M15").Value = -3998.1413
$ws.Range("H70").Value = 9864000
$ws.Range("I70").Value = 41917492
$ws.Range("J70").Value = 1386.4615
$ws.Range("K70").Value = 125752476
$ws.Range("L70").Value = 4159.3845
$ws.Range("M70").Value = -125752206
$ws.Range("N70").Value = -4699.3845
$ws.Range("H73").Value = 9864000
$ws.Range("I73").Value = 41917492
$ws.Range("J73").Value = 1386.4615
$ws.Range("K73").Value = 125752476
$ws.Range("L73").Value = 4159.3845
$ws.Range("M73").Value = -125751540
$ws.Range("N73").Value = -6031.3845
$ws.Range("H100").Value = 2169.125
$ws.Range("I100").Value = 1213.125
$ws.Range("K100").Value = 1213.125
$ws.Range("M100").Value = -672.125
$ws.Range("H113").Value = 2888.7778
$ws.Range("J113").Value = 3333.3333
$ws.Range("L113").Value = 3333.3333
$ws.Range("N113").Value = -9841.3333
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H53").Value = 25000
$ws.Range("J53").Value = 25000
$ws.Range("L53").Value = 25000
$ws.Range("N53").Value = -26364
$ws.Range("H61").Value = 1904.4166
$ws.Range("I61").Value = 1605.3
$ws.Range("K61").Value = 1605.3
$ws.Range("M61").Value = -1393.3
$ws.Range("H96").Value = 155499.75
$ws.Range("J96").Value = 155499.75
$ws.Range("L96").Value = 155499.75
$ws.Range("N96").Value = -160991.75
$ws.Range("H101").Value = 50951.5
$ws.Range("J101").Value = 50951.5
$ws.Range("L101").Value = 50951.5
$ws.Range("N101").Value = -57441.5
$ws.Range("H132").Value = 2327
$ws.Range("I132").Value = 1733.9333
$ws.Range("J132").Value = 3439
$ws.Range("K132").Value = 5201.7999
$ws.Range("L132").Value = 10317
$ws.Range("M132").Value = -2671.7999
$ws.Range("N132").Value = -15377
$ws.Range("H136").Value = 1904.4166
$ws.Range("I136").Value = 1605.3
$ws.Range("K136").Value = 4815.9
$ws.Range("M136").Value = -2265.9
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1875.125
$ws.Range("I16").Value = 1842.2
$ws.Range("J16").Value = 1930
$ws.Range("K16").Value = 1842.2
$ws.Range("L16").Value = 1930
$ws.Range("M16").Value = -1555.2
$ws.Range("N16").Value = -2504
$ws.Range("H113").Value = 1875.125
$ws.Range("I113").Value = 1842.2
$ws.Range("J113").Value = 1930
$ws.Range("K113").Value = 1842.2
$ws.Range("L113").Value = 1930
$ws.Range("M113").Value = 327.8
$ws.Range("N113").Value = -6270
$ws.Range("H132").Value = 2356.2144
$ws.Range("I132").Value = 1639.909
$ws.Range("J132").Value = 4982.6665
$ws.Range("K132").Value = 4919.727000000001
$ws.Range("L132").Value = 14947.9995
$ws.Range("M132").Value = -2389.727000000001
$ws.Range("N132").Value = -20007.9995
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1455.3
$ws.Range("I68").Value = 1375
$ws.Range("J68").Value = 1475.375
$ws.Range("K68").Value = 4125
$ws.Range("L68").Value = 4426.125
$ws.Range("M68").Value = -3314
$ws.Range("N68").Value = -6048.125
$ws.Range("H71").Value = 1455.3
$ws.Range("I71").Value = 1375
$ws.Range("J71").Value = 1475.375
$ws.Range("K71").Value = 12375
$ws.Range("L71").Value = 13278.375
$ws.Range("M71").Value = -8319
$ws.Range("N71").Value = -21390.375
$ws.Range("H137").Value = 41669176
$ws.Range("I137").Value = 1200
$ws.Range("J137").Value = 47621744
$ws.Range("K137").Value = 3600
$ws.Range("L137").Value = 142865232
$ws.Range("M137").Value = 1500
$ws.Range("N137").Value = -142875432
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2458.6667
$ws.Range("I107").Value = 5451
$ws.Range("J107").Value = 962.5
$ws.Range("K107").Value = 5451
$ws.Range("L107").Value = 962.5
$ws.Range("M107").Value = -3531
$ws.Range("N107").Value = -4802.5
$ws.Range("H122").Value = 2985.182
$ws.Range("I122").Value = 3117.4443
$ws.Range("J122").Value = 2390
$ws.Range("K122").Value = 9352.332900000001
$ws.Range("L122").Value = 7170
$ws.Range("M122").Value = -6902.332900000001
$ws.Range("N122").Value = -12070
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2406.5334
$ws.Range("I68").Value = 1399.4445
$ws.Range("J68").Value = 3917.1667
$ws.Range("K68").Value = 1399.4445
$ws.Range("L68").Value = 3917.1667
$ws.Range("M68").Value = -650.4445000000001
$ws.Range("N68").Value = -5415.1667
$ws.Range("H71").Value = 2406.5334
$ws.Range("I71").Value = 1399.4445
$ws.Range("J71").Value = 3917.1667
$ws.Range("K71").Value = 6997.2225
$ws.Range("L71").Value = 19585.8335
$ws.Range("M71").Value = -3253.2225
$ws.Range("N71").Value = -27073.8335
$ws.Range("H100").Value = 8225.375
$ws.Range("I100").Value = 11460.6
$ws.Range("J100").Value = 2833.3333
$ws.Range("K100").Value = 11460.6
$ws.Range("L100").Value = 2833.3333
$ws.Range("M100").Value = -10919.6
$ws.Range("N100").Value = -3915.3333
$ws.Range("H106").Value = 42999
$ws.Range("J106").Value = 42999
$ws.Range("L106").Value = 42999
$ws.Range("N106").Value = -45523
$ws.Range("H136").Value = 2708.9614
$ws.Range("I136").Value = 2342.6365
$ws.Range("J136").Value = 4723.75
$ws.Range("K136").Value = 7027.9095
$ws.Range("L136").Value = 14171.25
$ws.Range("M136").Value = -4477.9095
$ws.Range("N136").Value = -19271.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 29979.75
$ws.Range("J31").Value = 29979.75
$ws.Range("L31").Value = 29979.75
$ws.Range("N31").Value = -30675.75
$ws.Range("H132").Value = 1268.125
$ws.Range("I132").Value = 978.0909
$ws.Range("J132").Value = 1906.2
$ws.Range("K132").Value = 2934.2727
$ws.Range("L132").Value = 5718.6
$ws.Range("M132").Value = -404.2727
$ws.Range("N132").Value = -10778.6
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 918.75
$ws.Range("I105").Value = 918.75
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 918.75
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 828.25
$ws.Range("N105").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 131.33333
$ws.Range("I55").Value = 150
$ws.Range("J55").Value = 94
$ws.Range("K55").Value = 150
$ws.Range("L55").Value = 94
$ws.Range("M55").Value = 23
$ws.Range("N55").Value = -440
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
